$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "218.62") need the
# column pre-set to Text format, otherwise Excel silently reinterprets the
# assignment as a numeric value (losing formatting like trailing zeros,
# e.g. "62.60" -> 62.6) instead of keeping the literal price string.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.193.01'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '1.647.34'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '218.62'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").Value = '0.256'
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("D9").Value = '0.0627'
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("D10").Value = '20.08'
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '1.878.98'
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").Value = '1.623.44'
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").Value = '4.13'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("D15").Value = '0.538'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = '67.49'
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").Value = '27.179.54'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = '0.0₃0741'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '219.95'
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = '6.76'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '4.44'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '2.52'
$ws.Range("E23").Value = '  +4.00%  '
$ws.Range("D24").Value = '9.22'
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("D25").Value = '148.29'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '7.39'
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("E28").Value = '  -0.54%  '
$ws.Range("D29").Value = '15.81'
$ws.Range("E29").Value = '  -1.36%  '
$ws.Range("D30").Value = '0.0506'
$ws.Range("E30").Value = '  -1.79%  '
$ws.Range("E31").Value = '  -1.51%  '
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("D34").Value = '1.58'
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("D35").Value = '1.271.74'
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = '2.45'
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("D38").Value = '0.543'
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("D39").Value = '0.846'
$ws.Range("E39").Value = '  +1.20%  '
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("D41").Value = '0.811'
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("E42").Value = '  +3.36%  '
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("D44").Value = '1.789.41'
$ws.Range("E44").Value = '  -1.10%  '
$ws.Range("D45").Value = '62.60'
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").Value = '92.26'
$ws.Range("E46").Value = '  -0.51%  '
$ws.Range("D47").Value = '1.59'
$ws.Range("E47").Value = '  -2.00%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0514'
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '7.72'
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.0972'
$ws.Range("E50").Value = '  -1.15%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.406'
$ws.Range("E51").Value = '  -0.40%  '
